# Template News.pptx - update the "Flusso principale" description in the
# use-case table on slide 4: the system now fetches news from "InfoBlu"
# instead of simply "refreshing current news".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item("Tabella 4")
$tbl = $sh.Table
$cell = $tbl.Cell(1, 2)
$tr = $cell.Shape.TextFrame.TextRange
$paragraphs = $tr.Paragraphs()

# Second bullet of the first row: "Il sistema aggiornerà le news attuali"
# becomes "Il sistema reperisce le news da InfoBlu".
$paragraphs.Item(2).Text = "Il sistema reperisce le news da InfoBlu"
